$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Q4: new "2020" year header cell ---
# Reuse the exact formatting (font/border/number format) of the existing
# year-header cells (D4:H4, L4) via copy/paste-special, then nudge the
# vertical alignment from "center" to "top" to match the new style.
$ws.Range("D4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020
$ws.Range("Q4").VerticalAlignment = -4160

# --- Q5: new "21.8" data cell, matching the rest of row 5's number cells ---
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 21.8

$excel.CutCopyMode = $false

# --- selection moves to Q9 ---
$ws.Range("Q9").Select()

Write-Output "done"
